$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 4 with the latest stock data snapshot, matching the
# string/numeric typing of the existing rows (A-F as text, G-N as numbers).
$ws.Cells.Item(4, 1).Value = "Tue Oct 31 00:50:03 2023"
$ws.Cells.Item(4, 2).Value = "HKHSCEI"
$ws.Cells.Item(4, 3).Value = "国企指数"
$ws.Cells.Item(4, 4).Value = "未开盘"
$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = "5960.33"
$ws.Cells.Item(4, 5).Style = "Normal"
$ws.Cells.Item(4, 6).Value = "-18.68  -0.31%"
$ws.Cells.Item(4, 7).Value = 5972.03
$ws.Cells.Item(4, 8).Value = 5915.54
$ws.Cells.Item(4, 9).Value = 7773.61
$ws.Cells.Item(4, 10).Value = 0.012
$ws.Cells.Item(4, 11).Value = 5897.6
$ws.Cells.Item(4, 12).Value = 5979.01
$ws.Cells.Item(4, 13).Value = 4919.03
$ws.Cells.Item(4, 14).Value = 39705000000
